$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (B3/C3): was "1_DB_Script"/"Anson Antony", becomes "Dummy_Data"/"Soumya Raj"
$ws.Range("B3").Value = "Dummy_Data"
$ws.Range("C3").Value = "Soumya Raj"

# Update row 4 (A4 stays 3; B4/C4 now hold what used to be in row 3: "1_DB_Script"/"Anson Antony")
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "1_DB_Script"
$ws.Range("C4").Value = "Anson Antony"

# Remove now-unused rows 5 and 6 (previously held SL NO 4 and 5 with no other data)
$ws.Range("A5:C6").ClearContents()

# Update selection to match the target state
$ws.Range("B4").Select()
